$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) and Volume (E) columns so numeric-looking
# strings (e.g. "213.22", "1.00", "  +0.24%  ") are preserved verbatim as text
# instead of being auto-converted to Excel numbers.
$textCells = @(
  "D2",
  "E2",
  "D3",
  "E3",
  "E4",
  "D5",
  "E5",
  "D6",
  "E6",
  "E7",
  "E8",
  "E9",
  "D10",
  "E10",
  "E11",
  "D12",
  "E12",
  "D13",
  "E13",
  "E14",
  "E15",
  "D16",
  "E16",
  "D17",
  "E17",
  "E18",
  "D19",
  "E19",
  "E20",
  "E21",
  "E22",
  "D23",
  "E23",
  "D24",
  "E24",
  "D25",
  "E25",
  "D26",
  "E26",
  "E27",
  "E28",
  "D29",
  "E29",
  "D30",
  "E30",
  "E31",
  "E32",
  "D33",
  "E33",
  "D34",
  "E34",
  "D35",
  "E35",
  "E37",
  "E38",
  "E39",
  "D40",
  "E40",
  "E41",
  "D42",
  "E42",
  "E43",
  "D44",
  "E44",
  "D45",
  "E45",
  "E46",
  "D47",
  "E47",
  "D48",
  "E48",
  "D49",
  "E49",
  "D50",
  "E50",
  "D51",
  "E51",
)
foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.520.97"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "1.627.75"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "213.22"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "0.505"
$ws.Range("E6").Value = "  +2.17%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "18.79"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").Value = "1.853.86"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "1.631.87"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").Value = "65.13"
$ws.Range("E16").Value = "  +3.38%  "
$ws.Range("D17").Value = "26.547.45"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "214.76"
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D23").Value = "9.32"
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("D24").Value = "2.14"
$ws.Range("E24").Value = "  +12.48%  "
$ws.Range("D25").Value = "147.62"
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("D29").Value = "15.58"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").Value = "0.0513"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("D33").Value = "2.97"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "1.50"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "1.239.46"
$ws.Range("E35").Value = "  +6.00%  "
$ws.Range("E37").Value = "  +4.42%  "
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("D40").Value = "0.794"
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("E41").Value = "  -2.21%  "
$ws.Range("D42").Value = "0.798"
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").Value = "1.763.86"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").Value = "93.13"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("E46").Value = "  +2.57%  "
$ws.Range("D47").Value = "54.91"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0510"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0959"
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.407"
$ws.Range("E51").Value = "  -0.53%  "
